# Apply cryptos list update (Tue Sep 19 20:12:27 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.188.99'
$ws.Range('D3').Value = '1.643.76'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''217.07'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = '''0.514'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('D10').Value = '''19.91'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').Value = '1.876.24'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '1.647.53'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').Value = '''0.542'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('D16').Value = '''67.33'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').Value = '27.184.09'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').Value = '''218.83'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '''6.83'
$ws.Range('E21').Value = '  +2.88%  '
$ws.Range('D22').Value = '''2.55'
$ws.Range('E22').Value = '  +4.73%  '
$ws.Range('D23').Value = '''4.39'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '''9.17'
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = '''147.80'
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('D26').Value = '''7.55'
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = '''15.74'
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').Value = '''1.19'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').Value = '''3.03'
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('D34').Value = '''1.57'
$ws.Range('E34').Value = '  +0.87%  '
$ws.Range('D35').Value = '1.259.50'
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').Value = '''0.0177'
$ws.Range('E37').Value = '  +1.82%  '
$ws.Range('D38').Value = '''0.543'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').Value = '''0.846'
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D41').Value = '''0.807'
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('D42').Value = '''2.23'
$ws.Range('E42').Value = '  +6.48%  '
$ws.Range('D43').Value = '''5.37'
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').Value = '1.785.73'
$ws.Range('D45').Value = '''61.63'
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('D46').Value = '''91.63'
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = '''1.60'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.66'
$ws.Range('E50').Value = '  +1.14%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.0971'
$ws.Range('E51').Value = '  -0.33%  '
